# Apply "Derniere mise du journal de bord" edit:
# 1) Append a new run to the end of the last existing paragraph
#    ("Vendredi 5 decembre 2014" entry).
# 2) Insert a large block of new dated journal paragraphs
#    (Samedi 6 decembre 2014 .. Lundi 22 decembre 2014) plus a
#    trailing empty paragraph.
# 3) Relocate the "_GoBack" bookmark so that it now sits inside
#    the "Lundi 15 decembre 2014" entry, splitting "texture"
#    into "textur" + "e du missile".

$d = $word.ActiveDocument

# The existing _GoBack bookmark sits right at the end of the last
# paragraph; remove it first so new content is not pushed past it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 1: extend the last paragraph with the new sentence ---
$lastPara = $d.Paragraphs.Last
$lastParaRange = $lastPara.Range
$textEnd = $lastParaRange.End - 1
$insertionPoint = $d.Range($textEnd, $textEnd)
$newSentence = " Réglage de bug avec le joueur et il n’y a plus de changement de tour à chaque fois qu’on appuie sur une flèche."
$insertionPoint.InsertAfter($newSentence)
$newSentenceRange = $d.Range($textEnd, $textEnd + $newSentence.Length)
$newSentenceRange.Font.Size = 12

# --- Step 2: append the whole block of new journal paragraphs ---
$endPara = $d.Paragraphs.Last
$endRange = $endPara.Range
$appendPoint = $d.Range($endRange.End, $endRange.End)
$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Samedi 6 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">J'ai travaillé sur le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>jetpack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> et sur le m</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>enu pause. J'ai mis la physique</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> des pa</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">cks dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CGame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">, pour mettre la </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>régression aux packs.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Lundi 8 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai travaillé sur les collisions avec le joueur et le missile.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Mardi 9 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai travaillé sur le missile, sa physique et ses collisions</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>, mais le missile n’est pas encore fini</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Mercredi 10 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai terminé le missile, mais il reste à perfectionner les collisions avec le joueur</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Jeudi 11 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai ajouté les collisions avec le joueur, mais ce n'est pas fini.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Samedi 13 décembre 2014 :</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> J'ai continué le missile, mais il reste à régler quelque bug avec les collisions</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Dimanche 14 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Le missile est maintenant terminé, les dommages sont maintenant parfait</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Lundi 15 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai perfectionné les dommages et j'ai fait affich</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>er</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> les </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>rects</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>player</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> et des ex</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">plosions. J'ai aussi </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>changé</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> la </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>texture du missile</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Mardi 16 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>La tran</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>parence se détecte automatiquement pour chaque surface. J'ai travaillé sur le bazo</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">oka, il s'affiche maintenant et </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>les mines explosent avant de disparaître.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Mercredi 17 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Les dommages sont encore améliorés, le programme plante plus lors de la sup</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>p</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>ression d'un personnage et les mines e</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">xplosent </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">en chaîne. J'ai aussi travaillé sur la physique du </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>jetpack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">. J'ai aussi régler un bug créé avec des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>merges</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Jeudi 18 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Le menu principal à un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>bacground</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">. J'ai fini le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>jectpack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>, il reste maintenant à gérer les collisions avec le joueur.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Vendredi 19 décembre 2014 </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">J'ai fait les </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>destructeurs</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>, mais il</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> ont disparus... Les dommages s'affichent maintenant. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Les </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>tools</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> se réinitialisent au </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>début de chaque tour (barre de puissance et vecteur de vitesse).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Samedi 20 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">J'ai réglé des problèmes avec les dommages, il n'y a plus de zéro qui </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>affichent</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>. Ajout du dommage e</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">n chute libre (enfin...). </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>sprite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> d'explosion s'affiche tou</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> le temps. Le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>sprite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> se remet au départ quand il n'est pas en bouc</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">le, sinon lors de sa prochaine </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">utilisation, il ne s'affichera pas. J'ai réglé quelque bug avec le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>healthpack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>. Il n'y a plus zéro de d</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ommage quand un missile atteint </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">un joueur. On réinitialise le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>rect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> source de la barre de puissance pour qu'il s'affiche correctement lo</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">rs de sa prochaine utilisation. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">On gère la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>toolbarre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CGame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> au lieu de dans la main.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Dimanche 21 décembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">J'ai travaillé sur la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>mélée</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>, il est maintenant un projectile.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Lundi 22 décembre 2014 :</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>J'ai fini l</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> mêlée. J'ai perfectionné le projet, tout en complétant les commentaires.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr></w:p>
'@
[void]$appendPoint.InsertXML($newParagraphsXml)

# --- Step 3: put the _GoBack bookmark back, inside "texture du missile" ---
$searchRange = $d.Content
$found = $searchRange.Find.Execute("texture du missile", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPos = $searchRange.Start + "textur".Length
    $bookmarkRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

Write-Output "Journal de bord updated."
